$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 12:55"

# Espana (row 6): F6 5872 -> 6092
$ws.Range("F6").Value = 6092

# Alemania (row 8)
$ws.Range("B8").Value = 78115
$ws.Range("C8").Value = 134
$ws.Range("E8").Value = 57996
$ws.Range("G8").Value = 13
$ws.Range("H8").Value = 944

# Row 69/70: Barein overtakes Hungria
$ws.Range("A69").Value = "Barein"
$ws.Range("B69").Value = 635
$ws.Range("C69").Value = 66
$ws.Range("D69").Value = 341
$ws.Range("E69").Value = 290
$ws.Range("F69").Value = 3
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 4

$ws.Range("A70").Value = "Hungria"
$ws.Range("B70").Value = 585
$ws.Range("C70").Value = 60
$ws.Range("D70").Value = 42
$ws.Range("E70").Value = 522
$ws.Range("F70").Value = 17
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 21

# Rows 98-102: Malta overtakes Ghana, Senegal, Costa de Marfil, Uzbekistan
$ws.Range("A98").Value = "Malta"
$ws.Range("B98").Value = 196
$ws.Range("C98").Value = 8
$ws.Range("D98").Value = 2
$ws.Range("E98").Value = 194
$ws.Range("F98").Value = 2
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0

$ws.Range("A99").Value = "Ghana"
$ws.Range("B99").Value = 195
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 31
$ws.Range("E99").Value = 159
$ws.Range("F99").Value = 1
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 5

$ws.Range("A100").Value = "Senegal"
$ws.Range("B100").Value = 195
$ws.Range("C100").Value = 5
$ws.Range("D100").Value = 55
$ws.Range("E100").Value = 139
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 1

$ws.Range("A101").Value = "Costa de Marfil"
$ws.Range("B101").Value = 190
$ws.Range("C101").Value = 0
$ws.Range("D101").Value = 9
$ws.Range("E101").Value = 180
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 1

$ws.Range("A102").Value = "Uzbekistan"
$ws.Range("B102").Value = 190
$ws.Range("C102").Value = 9
$ws.Range("D102").Value = 12
$ws.Range("E102").Value = 176
$ws.Range("F102").Value = 8
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 2
